$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The bioSampleNumber column (C) collided with another workbook (0648); bump
# every value in C2:C19 by 17 so the keys become unique (1..18 -> 18..35).
for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 17
}

# Touch C20 so a (blank) trailing row 20 exists and the sheet's used range
# grows to include it, then give every data row (2-20) the new row height.
$ws.Cells.Item(20, 3).NumberFormat = "General"
for ($r = 2; $r -le 20; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Move the active selection to E8.
$ws.Range("E8").Select() | Out-Null
